# Auto-generated edit script for horarios-141-2026-01-16.xlsx
# Updates scrape data (new scrape at 06:02:16) across the three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912" (sheet1): header + full replacement of data rows 31-64
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:02:16"
$ws1.Range("A3").Value = "Total filas: 59"

$sheet1Rows = @(
    @("06:02:16","06:05","16_SANTA ANA",3,"LP1912"),
    @("04:18:06","06:09","16_SANTA ANA",111,"LP1912"),
    @("04:40:32","06:11","215A_EL PATO",91,"LP1912"),
    @("04:18:06","06:12","215A_EL PATO",114,"LP1912"),
    @("04:18:06","06:14","225_HARAS DEL SUR",116,"LP1912"),
    @("04:40:32","06:21","26_HERNANDEZ",101,"LP1912"),
    @("04:40:32","06:27","23_HERNANDEZ",107,"LP1912"),
    @("04:40:32","06:29","86_EST CHICA-ESC AGRARIA",109,"LP1912"),
    @("06:02:16","06:30","86_EST CHICA-ESC AGRARIA",28,"LP1912"),
    @("04:40:32","06:31","16_SANTA ANA",111,"LP1912"),
    @("04:53:50","06:44","225_C ROCA-H SUR",111,"LP1912"),
    @("04:53:50","06:46","215C_EL PATO",113,"LP1912"),
    @("05:18:23","06:58","10_OLMOS",100,"LP1912"),
    @("05:18:23","06:59","14_ABASTO",101,"LP1912"),
    @("06:02:16","07:00","14_ABASTO",58,"LP1912"),
    @("05:47:32","07:04","23_HERNANDEZ",77,"LP1912"),
    @("05:18:23","07:05","15_ABASTO",107,"LP1912"),
    @("06:02:16","07:05","23_HERNANDEZ",63,"LP1912"),
    @("05:18:23","07:07","225_GOMEZ",109,"LP1912"),
    @("05:18:23","07:11","215A_EL PATO",113,"LP1912"),
    @("05:18:23","07:15","11_ETCHEVERRY",117,"LP1912"),
    @("06:02:16","07:16","11_ETCHEVERRY",74,"LP1912"),
    @("05:47:32","07:21","26_HERNANDEZ",94,"LP1912"),
    @("06:02:16","07:23","10_OLMOS",81,"LP1912"),
    @("05:47:32","07:27","10_OLMOS",100,"LP1912"),
    @("05:47:32","07:31","11_ETCHEVERRY",104,"LP1912"),
    @("05:47:32","07:31","16_SANTA ANA",104,"LP1912"),
    @("05:47:32","07:32","84_COLONIA URQUIZA-ESC 49",105,"LP1912"),
    @("06:02:16","07:32","11_ETCHEVERRY",90,"LP1912"),
    @("05:47:32","07:36","27_EL RETIRO",109,"LP1912"),
    @("06:02:16","07:37","27_EL RETIRO",95,"LP1912"),
    @("05:47:32","07:39","10_OLMOS",112,"LP1912"),
    @("06:02:16","07:48","14_ABASTO",106,"LP1912"),
    @("06:02:16","07:51","215D_EL PATO",109,"LP1912")
)

$startRow1 = 31
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $row = $sheet1Rows[$i]
    $r = $startRow1 + $i
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------
# Sheet "LP1912-215" (sheet2): header + append new row 17
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:02:16"
$ws2.Range("A3").Value = "Total filas: 12"

$ws2.Cells.Item(17, 1).Value = "06:02:16"
$ws2.Cells.Item(17, 2).Value = "07:51"
$ws2.Cells.Item(17, 3).Value = "215D_EL PATO"
$ws2.Cells.Item(17, 4).Value = 109
$ws2.Cells.Item(17, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "6203-6173" (sheet3): header timestamp only
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:02:16"
